$d = $word.ActiveDocument
$d.Content.Find.Execute("stop", $true, $true, $false, $false, $false,
                         $true, 1, $false, "", 2)
